$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 45 (everything from the old row 45 onward shifts down
# by one, matching the diff's renumbering of rows 45-135 -> 46-136).
$ws.Rows.Item(45).Insert() | Out-Null

# Seed the new row's formatting (value + style) from the row immediately
# above it (row 44: "대우목화" / 7378), which carries the same border/fill
# combination ("s=3" / white fill) that the new entry ends up using.
$ws.Range("A44:B44").Copy($ws.Range("A45:B45"))

# New apartment entry: ID 7306, name "현대".
$ws.Range("A45").Value = 7306
$ws.Range("B45").Value = "현대"

# The new entry's label uses a distinct font (Dotum 9pt) rather than the
# Arial 9pt used by the rest of column B.
$ws.Range("B45").Font.Name = "돋움"
$ws.Range("B45").Font.Size = 9

# Update the view: scrolled down with D44 as the active cell.
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("D44").Select() | Out-Null

# Page setup now explicit (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
